$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1095.3478
$ws.Range("I28").Value = 866.6667
$ws.Range("K28").Value = 866.6667
$ws.Range("M28").Value = -381.6667

$ws.Range("H33").Value = 12500463
$ws.Range("I33").Value = 12500463
$ws.Range("K33").Value = 12500463
$ws.Range("M33").Value = -12500234

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H103").Value = 364.14285
$ws.Range("J103").Value = 383
$ws.Range("L103").Value = 1149
$ws.Range("N103").Value = -2321

$ws.Range("H137").Value = 2740.9167
$ws.Range("I137").Value = 2435.8948
$ws.Range("K137").Value = 7307.6844
$ws.Range("M137").Value = -4757.6844

$ws.Range("H138").Value = 6669029.5
$ws.Range("I138").Value = 1018.11536
$ws.Range("J138").Value = 10207158
$ws.Range("K138").Value = 3054.34608
$ws.Range("L138").Value = 30621474
$ws.Range("M138").Value = 2085.65392
$ws.Range("N138").Value = -30631754

$ws.Range("H141").Value = 1115.5405
$ws.Range("I141").Value = 1134.3334
$ws.Range("K141").Value = 3403.0002
$ws.Range("M141").Value = 1776.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 14999
$ws.Range("J8").Value = 14999
$ws.Range("L8").Value = 14999
$ws.Range("N8").Value = -15287

$ws.Range("H32").Value = 5311.9395
$ws.Range("I32").Value = 2386.7378
$ws.Range("K32").Value = 2386.7378
$ws.Range("M32").Value = -2099.7378

$ws.Range("H45").Value = 5726.5938
$ws.Range("I45").Value = 6701.48
$ws.Range("J45").Value = 2244.8572
$ws.Range("K45").Value = 6701.48
$ws.Range("L45").Value = 2244.8572
$ws.Range("M45").Value = -6324.48
$ws.Range("N45").Value = -2998.8572

$ws.Range("H61").Value = 3011.2666
$ws.Range("I61").Value = 1955.5
$ws.Range("J61").Value = 5122.8
$ws.Range("K61").Value = 1955.5
$ws.Range("L61").Value = 5122.8
$ws.Range("M61").Value = -1743.5
$ws.Range("N61").Value = -5546.8

$ws.Range("H74").Value = 29679.62
$ws.Range("I74").Value = 32776.188
$ws.Range("K74").Value = 32776.188
$ws.Range("M74").Value = -31902.188

$ws.Range("H77").Value = 29679.62
$ws.Range("I77").Value = 32776.188
$ws.Range("K77").Value = 163880.94
$ws.Range("M77").Value = -159512.94

$ws.Range("H122").Value = 2344.2068
$ws.Range("I122").Value = 2298.76
$ws.Range("K122").Value = 6896.280000000001
$ws.Range("M122").Value = -4446.280000000001

$ws.Range("H136").Value = 3011.2666
$ws.Range("I136").Value = 1955.5
$ws.Range("J136").Value = 5122.8
$ws.Range("K136").Value = 5866.5
$ws.Range("L136").Value = 15368.4
$ws.Range("M136").Value = -3316.5
$ws.Range("N136").Value = -20468.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 3689.4443
$ws.Range("J11").Value = 4000
$ws.Range("L11").Value = 4000
$ws.Range("N11").Value = -4280

$ws.Range("H22").Value = 793.26666
$ws.Range("I22").Value = 630
$ws.Range("J22").Value = 1119.8
$ws.Range("K22").Value = 630
$ws.Range("L22").Value = 1119.8
$ws.Range("M22").Value = -280
$ws.Range("N22").Value = -1819.8

$ws.Range("H132").Value = 1696.8518
$ws.Range("I132").Value = 1647.32
$ws.Range("J132").Value = 2316
$ws.Range("K132").Value = 4941.96
$ws.Range("L132").Value = 6948
$ws.Range("M132").Value = -2411.96
$ws.Range("N132").Value = -12008

$ws.Range("H134").Value = 18455.756
$ws.Range("I134").Value = 7126.4194
$ws.Range("K134").Value = 21379.2582
$ws.Range("M134").Value = -18844.2582

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4561.1665
$ws.Range("I62").Value = 3341.75
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 10025.25
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = -9339.25
$ws.Range("N62").Value = -22372

$ws.Range("H65").Value = 4561.1665
$ws.Range("I65").Value = 3341.75
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 30075.75
$ws.Range("L65").Value = 63000
$ws.Range("M65").Value = -26643.75
$ws.Range("N65").Value = -69864

$ws.Range("H92").Value = 977.6875
$ws.Range("I92").Value = 598.5
$ws.Range("J92").Value = 1104.0834
$ws.Range("K92").Value = 1795.5
$ws.Range("L92").Value = 3312.2502
$ws.Range("M92").Value = -547.5
$ws.Range("N92").Value = -5808.2502

$ws.Range("H132").Value = 1337.9166
$ws.Range("I132").Value = 1247.258
$ws.Range("K132").Value = 11225.322
$ws.Range("M132").Value = -8695.322

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 9755.866
$ws.Range("J70").Value = 10009.5
$ws.Range("K70").Value = 9755.866
$ws.Range("L70").Value = 10009.5
$ws.Range("M70").Value = -9485.866
$ws.Range("N70").Value = -10549.5

$ws.Range("I73").Value = 9755.866
$ws.Range("J73").Value = 10009.5
$ws.Range("K73").Value = 9755.866
$ws.Range("L73").Value = 10009.5
$ws.Range("M73").Value = -8819.866
$ws.Range("N73").Value = -11881.5

$ws.Range("H102").Value = 43484316
$ws.Range("I102").Value = 1870
$ws.Range("K102").Value = 1870
$ws.Range("M102").Value = -248

$ws.Range("H107").Value = 544.8570999999999

$ws.Range("H122").Value = 2231.1555
$ws.Range("I122").Value = 2071.6875
$ws.Range("J122").Value = 2623.6924
$ws.Range("K122").Value = 6215.0625
$ws.Range("L122").Value = 7871.0772
$ws.Range("M122").Value = -3765.0625
$ws.Range("N122").Value = -12771.0772

$ws.Range("H132").Value = 4769.778
$ws.Range("I132").Value = 4797.067
$ws.Range("K132").Value = 14391.201
$ws.Range("M132").Value = -11861.201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3985.3333
$ws.Range("J7").Value = 6379.6
$ws.Range("L7").Value = 6379.6
$ws.Range("N7").Value = -6603.6

$ws.Range("H100").Value = 6206.5
$ws.Range("I100").Value = 2164.3333
$ws.Range("K100").Value = 2164.3333
$ws.Range("M100").Value = -1623.3333

$ws.Range("H126").Value = 3985.3333
$ws.Range("J126").Value = 6379.6
$ws.Range("L126").Value = 19138.8
$ws.Range("N126").Value = -24078.8

$ws.Range("H132").Value = 5207.25
$ws.Range("I132").Value = 4626.3335
$ws.Range("K132").Value = 13879.0005
$ws.Range("M132").Value = -11349.0005

$ws.Range("H136").Value = 2764.1345
$ws.Range("I136").Value = 2310
$ws.Range("J136").Value = 5261.875
$ws.Range("K136").Value = 6930
$ws.Range("L136").Value = 15785.625
$ws.Range("M136").Value = -4380
$ws.Range("N136").Value = -20885.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 19166.666
$ws.Range("J64").Value = 19166.666
$ws.Range("L64").Value = 19166.666
$ws.Range("N64").Value = -19662.666

$ws.Range("H67").Value = 19166.666
$ws.Range("J67").Value = 19166.666
$ws.Range("L67").Value = 19166.666
$ws.Range("N67").Value = -20882.666

$ws.Range("H122").Value = 1309.2325
$ws.Range("I122").Value = 1140.7297
$ws.Range("K122").Value = 3422.189100000001
$ws.Range("M122").Value = -972.1891000000005

$ws.Range("H132").Value = 197376.52
$ws.Range("I132").Value = 1302.7046
$ws.Range("J132").Value = 1429840.4
$ws.Range("K132").Value = 3908.1138
$ws.Range("L132").Value = 4289521.199999999
$ws.Range("M132").Value = -1378.1138
$ws.Range("N132").Value = -4294581.199999999
